$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend column A formatting (date style) from the last existing row down
# through the new rows, matching the existing pattern (style index used for
# dates in column A).
$ws.Range("A343").Copy()
$ws.Range("A344:A357").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(344, 1).Value = 44418
$ws.Cells.Item(344, 2).Value = 8
$ws.Cells.Item(344, 3).Value = 34
$ws.Cells.Item(344, 4).Value = 84.41332737474552
$ws.Cells.Item(345, 1).Value = 44419
$ws.Cells.Item(345, 2).Value = 2
$ws.Cells.Item(345, 3).Value = 35
$ws.Cells.Item(345, 4).Value = 86.89607229753216
$ws.Cells.Item(346, 1).Value = 44420
$ws.Cells.Item(346, 2).Value = 5
$ws.Cells.Item(346, 3).Value = 32
$ws.Cells.Item(346, 4).Value = 79.44783752917226
$ws.Cells.Item(347, 1).Value = 44421
$ws.Cells.Item(347, 2).Value = 9
$ws.Cells.Item(347, 3).Value = 37
$ws.Cells.Item(347, 4).Value = 91.86156214310542
$ws.Cells.Item(348, 1).Value = 44422
$ws.Cells.Item(348, 2).Value = 4
$ws.Cells.Item(348, 3).Value = 37
$ws.Cells.Item(348, 4).Value = 91.86156214310542
$ws.Cells.Item(349, 1).Value = 44423
$ws.Cells.Item(349, 2).Value = 6
$ws.Cells.Item(349, 3).Value = 40
$ws.Cells.Item(349, 4).Value = 99.30979691146531
$ws.Cells.Item(350, 1).Value = 44424
$ws.Cells.Item(350, 2).Value = 6
$ws.Cells.Item(350, 3).Value = 40
$ws.Cells.Item(350, 4).Value = 99.30979691146531
$ws.Cells.Item(351, 1).Value = 44425
$ws.Cells.Item(351, 2).Value = 4
$ws.Cells.Item(351, 3).Value = 36
$ws.Cells.Item(351, 4).Value = 89.37881722031878
$ws.Cells.Item(352, 1).Value = 44426
$ws.Cells.Item(352, 2).Value = 3
$ws.Cells.Item(352, 3).Value = 37
$ws.Cells.Item(352, 4).Value = 91.86156214310542
$ws.Cells.Item(353, 1).Value = 44427
$ws.Cells.Item(353, 2).Value = 4
$ws.Cells.Item(353, 3).Value = 36
$ws.Cells.Item(353, 4).Value = 89.37881722031878
$ws.Cells.Item(354, 1).Value = 44428
$ws.Cells.Item(354, 2).Value = 4
$ws.Cells.Item(354, 3).Value = 31
$ws.Cells.Item(354, 4).Value = 76.96509260638562
$ws.Cells.Item(355, 1).Value = 44429
$ws.Cells.Item(355, 2).Value = 6
$ws.Cells.Item(355, 3).Value = 33
$ws.Cells.Item(355, 4).Value = 81.93058245195888
$ws.Cells.Item(356, 1).Value = 44430
$ws.Cells.Item(356, 2).Value = 12
$ws.Cells.Item(356, 3).Value = 39
$ws.Cells.Item(356, 4).Value = 96.82705198867869
$ws.Cells.Item(357, 1).Value = 44431
$ws.Cells.Item(357, 2).Value = 0
$ws.Cells.Item(357, 3).Value = 33
$ws.Cells.Item(357, 4).Value = 81.93058245195888
